$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-free direct cell updates, preserving text (string) representation
# for numeric-looking values by forcing Text number format before assignment.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '318.14'

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '3.78%'

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '39.68'

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '1.96%'

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.139'

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '0.83%'

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.08211'

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '1.95%'

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '2.147'

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '11.37%'

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '8.306'

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '4.15%'

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9319'

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '0.04%'

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1414'

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-2.73%'

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1981'

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '2.76%'

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.09055'

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '0.08%'

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03477'

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-0.77%'

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09801'

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '0.29%'

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001402'

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '0.62%'

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.006147'

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '2.95%'

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.684'

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.290'

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '2.35%'

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.311'

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-2.89%'

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '0.81%'

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.1294'

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-2.46%'

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.900'

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '2.47%'

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-2.30%'

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04327'

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-1.25%'

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001226'

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-0.93%'

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004768'

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '11.55%'

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0001300'

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-0.11%'

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0003996'

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '-10.15%'

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02214'

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '8.72%'

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.05220'

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '3.81%'

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007501'

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '0.61%'

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.009680'

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-4.26%'

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '2.43%'

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.002114'

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-1.31%'

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.009849'

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '8.69%'

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006588'

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '6.55%'

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-0.17%'

$ws.Range("B48").Value = 'BOLO'

$ws.Range("C48").Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.002763'

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '-1.38%'

$ws.Range("B49").Value = 'CoinbaseStockToken'

$ws.Range("C49").Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.001200'

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-25.07%'

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.17%'

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '-0.17%'
